$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-106 down to 66-107
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the data for the inserted record
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44567
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 100112031
$ws.Range("G65").Value = "Poroto verde"
$ws.Range("H65").Value = "Brío"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 2000
$ws.Range("L65").Value = 2000
$ws.Range("M65").Value = 2000
$ws.Range("N65").Value = "$/kilo"
$ws.Range("O65").Value = "Región de La Araucanía"
$ws.Range("P65").Value = 2000
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
